$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
# continuing the "Modena" series through 27/05/2021.
$data = @(
    @(44330, 34, 204, 107.9290842428828),
    @(44331, 30, 207, 109.5162766582193),
    @(44332, 35, 213, 112.6906614888923),
    @(44333, 36, 210, 111.1034690735558),
    @(44334, 19, 217, 114.8069180426743),
    @(44335, 4, 208, 110.0453407966648),
    @(44336, 45, 203, 107.4000201044373),
    @(44337, 30, 199, 105.2837635506552),
    @(44338, 28, 197, 104.2256352737642),
    @(44339, 0, 162, 85.7083904281716),
    @(44340, 33, 159, 84.12119801283509),
    @(44341, 7, 147, 77.77242835148905),
    @(44342, 6, 149, 78.83055662838005),
    @(44343, 21, 125, 66.13301730568797)
)

# Last existing data row is 255 (date 13/05/2021); new rows start at 256.
$lastRow = 255
$r = $lastRow + 1
foreach ($row in $data) {
    # Clone column-A formatting (date style, border, bold, centered) from the
    # last existing row so the appended cell matches the rest of the series.
    $ws.Range("A$lastRow").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null

    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]

    $r++
}

$excel.CutCopyMode = 0
